$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9617735743522644
$ws.Range("B1").Value = 6.309230804443359
$ws.Range("C1").Value = 3.04810643196106
$ws.Range("D1").Value = 2.117587566375732
$ws.Range("E1").Value = 1.980319738388062
